# Add a "death_day" register column (C) to the "date" worksheet, mirroring
# the existing "birth_day" column (B): mostly dates, with one row holding a
# malformed "dd-dd" string instead of a real date (data-quality example).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("date")
$ws.Activate()

$ws.Range("C1").Value = "death_day"
$ws.Range("C2").Value = 45444
$ws.Range("C3").Value = 43954
$ws.Range("C4").Value = "12-56"
$ws.Range("C5").Value = 35033

# Format the real date cells the same way as the birth_day column.
$ws.Range("C2:C3").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("C5").NumberFormat = "yyyy\-mm\-dd"

$ws.Range("C5").Select()
